# Auto-generated edit script applying odds updates from the diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("K5").Value = 8
$ws.Range("N5").Value = 2.3
$ws.Range("O5").Value = 1.6
$ws.Range("G6").Value = 2.1
$ws.Range("H6").Value = 3.2
$ws.Range("I6").Value = 3.8
$ws.Range("P6").Value = 1.53
$ws.Range("Q6").Value = 2.38
$ws.Range("R6").Value = 2.05
$ws.Range("S6").Value = 1.7
$ws.Range("U6").Value = 9
$ws.Range("V6").Value = 9.5
$ws.Range("W6").Value = 19
$ws.Range("AF6").Value = 17
$ws.Range("G7").Value = 2.5
$ws.Range("R7").Value = 2
$ws.Range("S7").Value = 1.75
$ws.Range("W7").Value = 23
$ws.Range("Z7").Value = 7
$ws.Range("AB7").Value = 17
$ws.Range("AG7").Value = 12
$ws.Range("AH7").Value = 34
$ws.Range("N9").Value = 2.3
$ws.Range("O9").Value = 1.6
$ws.Range("G14").Value = 5.75
$ws.Range("H14").Value = 4.2
$ws.Range("I14").Value = 1.53
$ws.Range("P14").Value = 1.25
$ws.Range("Q14").Value = 3.75
$ws.Range("R14").Value = 1.53
$ws.Range("S14").Value = 2.38
$ws.Range("T14").Value = 21
$ws.Range("U14").Value = 34
$ws.Range("V14").Value = 17
$ws.Range("W14").Value = 51
$ws.Range("Y14").Value = 34
$ws.Range("AA14").Value = 8.5
$ws.Range("AC14").Value = 34
$ws.Range("AD14").Value = 101
$ws.Range("AF14").Value = 10
$ws.Range("AH14").Value = 13
$ws.Range("G15").Value = 2.45
$ws.Range("H15").Value = 3.25
$ws.Range("I15").Value = 2.88
$ws.Range("U15").Value = 13
$ws.Range("AI15").Value = 21
$ws.Range("G16").Value = 1.45
$ws.Range("H16").Value = 4.75
$ws.Range("I16").Value = 5.5
$ws.Range("L16").Value = 1.11
$ws.Range("M16").Value = 6.5
$ws.Range("N16").Value = 1.4
$ws.Range("O16").Value = 2.88
$ws.Range("U16").Value = 9.5
$ws.Range("V16").Value = 9
$ws.Range("W16").Value = 12
$ws.Range("AG16").Value = 17
$ws.Range("H17").Value = 5
$ws.Range("I17").Value = 9
$ws.Range("K17").Value = 15
$ws.Range("N17").Value = 1.53
$ws.Range("O17").Value = 2.4
$ws.Range("T17").Value = 8.5
$ws.Range("AA17").Value = 9.5
$ws.Range("N18").Value = 1.4
$ws.Range("O18").Value = 2.88
$ws.Range("J20").Value = 1.08
$ws.Range("K20").Value = 8
$ws.Range("G25").Value = 2.75
$ws.Range("H25").Value = 3.25
$ws.Range("I25").Value = 2.38
$ws.Range("J25").Value = 1.07
$ws.Range("K25").Value = 8.5
$ws.Range("N25").Value = 2.2
$ws.Range("O25").Value = 1.65
$ws.Range("Y25").Value = 41
$ws.Range("Z25").Value = 8.5
$ws.Range("G31").Value = 1.91
$ws.Range("I31").Value = 3.3
$ws.Range("R31").Value = 1.73
$ws.Range("S31").Value = 2
$ws.Range("X31").Value = 17
$ws.Range("AA31").Value = 6.5
$ws.Range("AD31").Value = 201
$ws.Range("AG31").Value = 12
$ws.Range("G34").Value = 1.44
$ws.Range("H34").Value = 4
$ws.Range("I34").Value = 6
$ws.Range("U34").Value = 7.5
$ws.Range("Y34").Value = 23
$ws.Range("AB34").Value = 17
$ws.Range("AD34").Value = 251
$ws.Range("AH34").Value = 67
$ws.Range("AI34").Value = 41
$ws.Range("AJ34").Value = 41
$ws.Range("K35").Value = 10
$ws.Range("N35").Value = 2.03
$ws.Range("O35").Value = 1.83
$ws.Range("P35").Value = 1.4
$ws.Range("Q35").Value = 2.75
$ws.Range("R35").Value = 1.75
$ws.Range("S35").Value = 2
$ws.Range("W35").Value = 21
$ws.Range("Z35").Value = 10
$ws.Range("AE35").Value = 9.5
$ws.Range("AH35").Value = 34
$ws.Range("H36").Value = 3.3
$ws.Range("N36").Value = 2
$ws.Range("O36").Value = 1.85
$ws.Range("R36").Value = 1.75
$ws.Range("S36").Value = 2
$ws.Range("T36").Value = 8.5
$ws.Range("U36").Value = 12
$ws.Range("X36").Value = 21
$ws.Range("Z36").Value = 10
$ws.Range("AC36").Value = 51
$ws.Range("AE36").Value = 9
$ws.Range("AF36").Value = 13
$ws.Range("AG36").Value = 11
$ws.Range("AI36").Value = 23
$ws.Range("X37").Value = 19
$ws.Range("Z37").Value = 8.5
$ws.Range("AE37").Value = 9.5
$ws.Range("AG37").Value = 12
$ws.Range("R38").Value = 1.64
$ws.Range("S38").Value = 2.17
$ws.Range("Y38").Value = 55
$ws.Range("L39").Value = 1.2
$ws.Range("M39").Value = 4.33
$ws.Range("N39").Value = 1.67
$ws.Range("O39").Value = 2.15
$ws.Range("G40").Value = 1.7
$ws.Range("H40").Value = 3.9
$ws.Range("K40").Value = 15
$ws.Range("L40").Value = 1.18
$ws.Range("M40").Value = 4.5
$ws.Range("N40").Value = 1.62
$ws.Range("O40").Value = 2.25
$ws.Range("P40").Value = 1.3
$ws.Range("Q40").Value = 3.4
$ws.Range("R40").Value = 1.62
$ws.Range("S40").Value = 2.2
$ws.Range("T40").Value = 9
$ws.Range("U40").Value = 9
$ws.Range("Z40").Value = 15
$ws.Range("AA40").Value = 7.5
$ws.Range("AD40").Value = 151
$ws.Range("AF40").Value = 26
$ws.Range("R43").Value = 1.53
$ws.Range("S43").Value = 2.38
$ws.Range("AF43").Value = 10
$ws.Range("H47").Value = 4.05
$ws.Range("I47").Value = 1.53
$ws.Range("J47").Value = 1.04
$ws.Range("K47").Value = 8.25
$ws.Range("L47").Value = 1.22
$ws.Range("M47").Value = 3.8
$ws.Range("N47").Value = 1.65
$ws.Range("O47").Value = 2.1
$ws.Range("P47").Value = 1.34
$ws.Range("Q47").Value = 3
$ws.Range("R47").Value = 1.75
$ws.Range("S47").Value = 1.95
$ws.Range("T47").Value = 17
$ws.Range("X47").Value = 50
$ws.Range("Z47").Value = 8.25
$ws.Range("AA47").Value = 7.9
$ws.Range("AB47").Value = 15.5
$ws.Range("AC47").Value = 65
$ws.Range("AD47").Value = 450
$ws.Range("AE47").Value = 7.5
$ws.Range("AF47").Value = 7.6
$ws.Range("AH47").Value = 11.25
$ws.Range("AI47").Value = 11.75
$ws.Range("AJ47").Value = 23
